# Update "想去人数" (column F) counts in both the "展览" and "全部类型" sheets
# to reflect the latest generated data (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 23
$ws1.Range("F3").Value = 74
$ws1.Range("F6").Value = 49
$ws1.Range("F7").Value = 2693
$ws1.Range("F9").Value = 276
$ws1.Range("F10").Value = 130
$ws1.Range("F11").Value = 10162
$ws1.Range("F12").Value = 75
$ws1.Range("F13").Value = 264
$ws1.Range("F16").Value = 11796
$ws1.Range("F17").Value = 12187

# --- Sheet "全部类型" (all types, combined listing) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 23
$ws4.Range("F3").Value = 74
$ws4.Range("F6").Value = 49
$ws4.Range("F7").Value = 2693
$ws4.Range("F10").Value = 276
$ws4.Range("F11").Value = 130
$ws4.Range("F12").Value = 10162
$ws4.Range("F13").Value = 75
$ws4.Range("F14").Value = 264
$ws4.Range("F17").Value = 11796
$ws4.Range("F18").Value = 12187
